# Commit: "updated the tab name"
# Rename the first worksheet (the test-cases sheet) from "LOGIN-TC-SHEET" to
# "TestCases-SHEET" and make it the active/selected tab (it was previously
# "VERSION-HISTORY" that held the selection).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LOGIN-TC-SHEET")
$ws.Name = "TestCases-SHEET"
$ws.Activate()
